$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "55.948.33"
$ws.Range("E2").Value = "  -1.59%  "

# Row 3
$ws.Range("D3").Value = "2.381.18"
$ws.Range("E3").Value = "  -5.10%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "478.96"
$ws.Range("E5").Value = "  -2.22%  "

# Row 6
$ws.Range("D6").Value = "146.86"
$ws.Range("E6").Value = "  -0.57%  "

# Row 7
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.22%  "

# Row 8
$ws.Range("E8").Value = "  -2.59%  "

# Row 9
$ws.Range("D9").Value = "2.384.05"
$ws.Range("E9").Value = "  -5.95%  "

# Row 10
$ws.Range("D10").Value = "0.0971"
$ws.Range("E10").Value = "  -1.12%  "

# Row 11
$ws.Range("D11").Value = "5.43"
$ws.Range("E11").Value = "  -5.93%  "

# Row 12
$ws.Range("D12").Value = "0.322"
$ws.Range("E12").Value = "  -3.50%  "

# Row 13
$ws.Range("E13").Value = "  +0.66%  "

# Row 14
$ws.Range("D14").Value = "2.795.80"
$ws.Range("E14").Value = "  -5.09%  "

# Row 15
$ws.Range("D15").Value = "55.994.23"
$ws.Range("E15").Value = "  -1.35%  "

# Row 16
$ws.Range("D16").Value = "20.28"
$ws.Range("E16").Value = "  -4.87%  "

# Row 17
$ws.Range("D17").Value = "0.0000131"
$ws.Range("E17").Value = "  -4.49%  "

# Row 18
$ws.Range("D18").Value = "2.398.10"
$ws.Range("E18").Value = "  -5.09%  "

# Row 19
$ws.Range("E19").Value = "  +0.15%  "

# Row 20
$ws.Range("D20").Value = "'314.20"
$ws.Range("E20").Value = "  -2.73%  "

# Row 21
$ws.Range("D21").Value = "9.66"
$ws.Range("E21").Value = "  -5.91%  "

# Row 22
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.20%  "

# Row 23
$ws.Range("D23").Value = "5.66"
$ws.Range("E23").Value = "  -3.14%  "

# Row 24
$ws.Range("D24").Value = "56.76"
$ws.Range("E24").Value = "  -3.86%  "

# Row 25
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.41%  "

# Row 26
$ws.Range("D26").Value = "0.394"
$ws.Range("E26").Value = "  -4.96%  "

# Row 27
$ws.Range("E27").Value = "  -6.68%  "

# Row 28
$ws.Range("D28").Value = "2.496.14"
$ws.Range("E28").Value = "  -4.79%  "

# Row 29
$ws.Range("D29").Value = "7.21"
$ws.Range("E29").Value = "  -5.72%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0765"
$ws.Range("E30").Value = "  -4.61%  "

# Row 31
$ws.Range("E31").Value = "  +0.03%  "

# Row 32
$ws.Range("D32").Value = "146.46"
$ws.Range("E32").Value = "  -1.94%  "

# Row 33
$ws.Range("D33").Value = "17.95"
$ws.Range("E33").Value = "  -2.16%  "

# Row 34
$ws.Range("E34").Value = "  -2.09%  "

# Row 35
$ws.Range("D35").Value = "5.03"
$ws.Range("E35").Value = "  -3.85%  "

# Row 36
$ws.Range("D36").Value = "'1.10"
$ws.Range("E36").Value = "  -4.84%  "

# Row 37
$ws.Range("D37").Value = "3.58"
$ws.Range("E37").Value = "  -5.06%  "

# Row 38
$ws.Range("D38").Value = "'0.830"
$ws.Range("E38").Value = "  -5.30%  "

# Row 39
$ws.Range("D39").Value = "33.38"
$ws.Range("E39").Value = "  -2.94%  "

# Row 40
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.23%  "

# Row 41
$ws.Range("D41").Value = "1.34"
$ws.Range("E41").Value = "  -0.58%  "

# Row 42
$ws.Range("D42").Value = "3.37"
$ws.Range("E42").Value = "  -5.54%  "

# Row 43
$ws.Range("D43").Value = "0.0535"
$ws.Range("E43").Value = "  -4.48%  "

# Row 44
$ws.Range("D44").Value = "0.0944"
$ws.Range("E44").Value = "  +3.00%  "

# Row 45
$ws.Range("E45").Value = "  -6.55%  "

# Row 46
$ws.Range("D46").Value = "10.21"
$ws.Range("E46").Value = "  +0.02%  "

# Row 47
$ws.Range("D47").Value = "253.45"
$ws.Range("E47").Value = "  -3.56%  "

# Row 48
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "0.0221"
$ws.Range("E48").Value = "  -3.57%  "

# Row 49
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "4.54"
$ws.Range("E49").Value = "  -5.96%  "

# Row 50
$ws.Range("D50").Value = "16.94"
$ws.Range("E50").Value = "  -4.72%  "

# Row 51
$ws.Range("D51").Value = "1.781.30"
$ws.Range("E51").Value = "  -7.94%  "
